# Add a new worksheet "Hex_Dec_Bin" with decimal/hex conversion example

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "Hex_Dec_Bin"

$newSheet.Range("A2").Value = "Décimal"
$newSheet.Range("B2").Value = -1

$newSheet.Range("A3").Value = "Hexadécimal"
$newSheet.Range("B3").Formula = "=DEC2HEX(B2)"

# Make the new sheet the active one
$newSheet.Activate()
$newSheet.Range("A4").Select()
